# EPBDS-10335 / EPBDS-8743 fix-up edit:
#  - Remove the (unused) Sheet2 / Sheet3 worksheets.
#  - Add a second "Test myRules" table (rows 27-30) below the existing one.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Drop the extra empty worksheets ---------------------------------
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# --- New "Test myRules" table (rows 27-30) -----------------------------
# Shared-string insertion order matches the original authoring order so the
# rebuilt sharedStrings.xml table lines up index-for-index with the target.
$ws.Range("B27").Value = "Test myRules"

$ws.Range("B28").Value = "aaa"
$ws.Range("C28").Value = "bbb"
$ws.Range("D28").Value = "_res_"

$ws.Range("D29").Value = "Result"
$ws.Range("B29").Value = "Arg1"
$ws.Range("C29").Value = "Arg2"

$ws.Range("B30").Value = "ccc"
$ws.Range("C30").Value = "yyy"

# --- Match the saved selection / active cell ---------------------------
$ws.Range("C30").Select() | Out-Null
